$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-7 with the new TPM-derived values (sending/target cluster
# reassignments plus recalculated expression/specificity numbers).
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Dlk1"
$ws.Cells.Item(2,3).Value = "Notch3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.101448
$ws.Cells.Item(2,8).Value = 6.304344
$ws.Cells.Item(2,9).Value = 0.5480341737688159
$ws.Cells.Item(2,10).Value = 0.5480341737688159
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 3.241087666666667
$ws.Cells.Item(2,14).Value = 9.723263
$ws.Cells.Item(2,15).Value = 0.02486257877280725
$ws.Cells.Item(2,16).Value = 0.02486257877280725
$ws.Cells.Item(2,17).Value = 6.810977194941333
$ws.Cells.Item(2,18).Value = 61.298794754472
$ws.Cells.Item(2,19).Value = 0.01362554281551752
$ws.Cells.Item(2,20).Value = 0.01362554281551752
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Dlk1"
$ws.Cells.Item(3,3).Value = "Notch3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.101448
$ws.Cells.Item(3,8).Value = 6.304344
$ws.Cells.Item(3,9).Value = 0.5480341737688159
$ws.Cells.Item(3,10).Value = 0.5480341737688159
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.732509
$ws.Cells.Item(3,14).Value = 8.197527
$ws.Cells.Item(3,15).Value = 0.02096124117795788
$ws.Cells.Item(3,16).Value = 0.02096124117795788
$ws.Cells.Item(3,17).Value = 5.742225573032
$ws.Cells.Item(3,18).Value = 51.68003015728799
$ws.Cells.Item(3,19).Value = 0.01148747649013103
$ws.Cells.Item(3,20).Value = 0.01148747649013103
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Dlk1"
$ws.Cells.Item(4,3).Value = "Notch3"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.101448
$ws.Cells.Item(4,8).Value = 6.304344
$ws.Cells.Item(4,9).Value = 0.5480341737688159
$ws.Cells.Item(4,10).Value = 0.5480341737688159
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 124.3864796666667
$ws.Cells.Item(4,14).Value = 373.159439
$ws.Cells.Item(4,15).Value = 0.9541761800492348
$ws.Cells.Item(4,16).Value = 0.9541761800492349
$ws.Cells.Item(4,17).Value = 261.3917189225573
$ws.Cells.Item(4,18).Value = 2352.525470303016
$ws.Cells.Item(4,19).Value = 0.5229211544631673
$ws.Cells.Item(4,20).Value = 0.5229211544631673
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Dlk1"
$ws.Cells.Item(5,3).Value = "Notch3"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.733072
$ws.Cells.Item(5,8).Value = 5.199216
$ws.Cells.Item(5,9).Value = 0.4519658262311841
$ws.Cells.Item(5,10).Value = 0.4519658262311841
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.241087666666667
$ws.Cells.Item(5,14).Value = 9.723263
$ws.Cells.Item(5,15).Value = 0.02486257877280725
$ws.Cells.Item(5,16).Value = 0.02486257877280725
$ws.Cells.Item(5,17).Value = 5.617038284645333
$ws.Cells.Item(5,18).Value = 50.553344561808
$ws.Cells.Item(5,19).Value = 0.01123703595728973
$ws.Cells.Item(5,20).Value = 0.01123703595728973
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Dlk1"
$ws.Cells.Item(6,3).Value = "Notch3"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.733072
$ws.Cells.Item(6,8).Value = 5.199216
$ws.Cells.Item(6,9).Value = 0.4519658262311841
$ws.Cells.Item(6,10).Value = 0.4519658262311841
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.732509
$ws.Cells.Item(6,14).Value = 8.197527
$ws.Cells.Item(6,15).Value = 0.02096124117795788
$ws.Cells.Item(6,16).Value = 0.02096124117795788
$ws.Cells.Item(6,17).Value = 4.735634837648
$ws.Cells.Item(6,18).Value = 42.62071353883199
$ws.Cells.Item(6,19).Value = 0.009473764687826854
$ws.Cells.Item(6,20).Value = 0.009473764687826854
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Dlk1"
$ws.Cells.Item(7,3).Value = "Notch3"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.733072
$ws.Cells.Item(7,8).Value = 5.199216
$ws.Cells.Item(7,9).Value = 0.4519658262311841
$ws.Cells.Item(7,10).Value = 0.4519658262311841
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 124.3864796666667
$ws.Cells.Item(7,14).Value = 373.159439
$ws.Cells.Item(7,15).Value = 0.9541761800492348
$ws.Cells.Item(7,16).Value = 0.9541761800492349
$ws.Cells.Item(7,17).Value = 215.5707250888694
$ws.Cells.Item(7,18).Value = 1940.136525799824
$ws.Cells.Item(7,19).Value = 0.4312550255860675
$ws.Cells.Item(7,20).Value = 0.4312550255860675

# Remove the former rows 8-10 (MuSCs sending-cluster block merged into
# the FAPs/MuSCs rows above under the new TPM data).
$ws.Rows("8:10").Delete()
